$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three species records that were on rows 2, 3 and 4 get reshuffled:
#   new row 2 <- old row 4 (Knärot / Goodyera repens record)
#   new row 3 <- old row 2 (Knärot / Goodyera repens record)
#   new row 4 <- old row 3 (Garnlav / Alectoria sarmentosa record, plus its comment)

# Row 2 <- old row 4
$ws.Range("A2").Value = 111790625
$ws.Range("B2").Value = 96348
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("Q2").Value = 489824.6884970492
$ws.Range("R2").Value = 6949020.70113107
$ws.Range("Z2").Value = "18:29"
$ws.Range("AB2").Value = "18:29"

# Row 3 <- old row 2
$ws.Range("A3").Value = 111792337
$ws.Range("B3").Value = 96348
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("Q3").Value = 489763.7116335144
$ws.Range("R3").Value = 6949091.647604217
$ws.Range("Z3").Value = "19:22"
$ws.Range("AB3").Value = "19:22"
$ws.Range("AC3").ClearContents()

# Row 4 <- old row 3 (including its public comment)
$ws.Range("A4").Value = 111790785
$ws.Range("B4").Value = 77515
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 489818.2822038208
$ws.Range("R4").Value = 6949032.207674611
$ws.Range("Z4").Value = "18:34"
$ws.Range("AB4").Value = "18:34"
$ws.Range("AC4").Value = "Många träd med mycket lav i området"
